# Apply data fixes for Team-Data/2010-11/5-5-2010-11.xlsx
# Commit: Fix Training Data Issue (#48) - stats shifted by one day
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the BF (Date) column to remain plain text so values like
# "2011-05-05" are not auto-converted into Excel date serials.
$ws.Range("BF2:BF31").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = 82
$ws.Cells.Item(2, 5).Value = 44
$ws.Cells.Item(2, 7).Value = 0.537
$ws.Cells.Item(2, 11).Value = 0.462
$ws.Cells.Item(2, 14).Value = 0.352
$ws.Cells.Item(2, 15).Value = 16.4
$ws.Cells.Item(2, 16).Value = 21.1
$ws.Cells.Item(2, 17).Value = 0.779
$ws.Cells.Item(2, 19).Value = 30
$ws.Cells.Item(2, 20).Value = 39.3
$ws.Cells.Item(2, 23).Value = 6.1
$ws.Cells.Item(2, 24).Value = 4.2
$ws.Cells.Item(2, 26).Value = 19
$ws.Cells.Item(2, 27).Value = 18.5
$ws.Cells.Item(2, 28).Value = 95
$ws.Cells.Item(2, 29).Value = -0.8
$ws.Cells.Item(2, 32).Value = 13
$ws.Cells.Item(2, 33).Value = 13
$ws.Cells.Item(2, 37).Value = 12
$ws.Cells.Item(2, 43).Value = 6
$ws.Cells.Item(2, 46).Value = 28
$ws.Cells.Item(2, 50).Value = 28
$ws.Cells.Item(2, 52).Value = 1
$ws.Cells.Item(2, 53).Value = 29
$ws.Cells.Item(2, 54).Value = 26
$ws.Cells.Item(2, 55).Value = 16
$ws.Cells.Item(2, 58).Value = "2011-05-05"

# Row 3
$ws.Cells.Item(3, 4).Value = 82
$ws.Cells.Item(3, 5).Value = 56
$ws.Cells.Item(3, 6).Value = 26
$ws.Cells.Item(3, 7).Value = 0.6830000000000001
$ws.Cells.Item(3, 8).Value = 48.2
$ws.Cells.Item(3, 9).Value = 36.9
$ws.Cells.Item(3, 10).Value = 75.8
$ws.Cells.Item(3, 11).Value = 0.486
$ws.Cells.Item(3, 14).Value = 0.365
$ws.Cells.Item(3, 17).Value = 0.77
$ws.Cells.Item(3, 22).Value = 14.6
$ws.Cells.Item(3, 23).Value = 8.199999999999999
$ws.Cells.Item(3, 26).Value = 20.5
$ws.Cells.Item(3, 28).Value = 96.5
$ws.Cells.Item(3, 29).Value = 5.4
$ws.Cells.Item(3, 30).Value = 1
$ws.Cells.Item(3, 32).Value = 6
$ws.Cells.Item(3, 42).Value = 22
$ws.Cells.Item(3, 43).Value = 11
$ws.Cells.Item(3, 47).Value = 4
$ws.Cells.Item(3, 48).Value = 22
$ws.Cells.Item(3, 55).Value = 6
$ws.Cells.Item(3, 58).Value = "2011-05-05"

# Row 4
$ws.Cells.Item(4, 4).Value = 82
$ws.Cells.Item(4, 6).Value = 48
$ws.Cells.Item(4, 7).Value = 0.415
$ws.Cells.Item(4, 9).Value = 35
$ws.Cells.Item(4, 10).Value = 77.59999999999999
$ws.Cells.Item(4, 11).Value = 0.451
$ws.Cells.Item(4, 14).Value = 0.327
$ws.Cells.Item(4, 17).Value = 0.756
$ws.Cells.Item(4, 18).Value = 10.3
$ws.Cells.Item(4, 20).Value = 40.1
$ws.Cells.Item(4, 25).Value = 6
$ws.Cells.Item(4, 28).Value = 93.3
$ws.Cells.Item(4, 29).Value = -4
$ws.Cells.Item(4, 34).Value = 18
$ws.Cells.Item(4, 43).Value = 21
$ws.Cells.Item(4, 44).Value = 21
$ws.Cells.Item(4, 48).Value = 21
$ws.Cells.Item(4, 52).Value = 8
$ws.Cells.Item(4, 53).Value = 13
$ws.Cells.Item(4, 58).Value = "2011-05-05"

# Row 5
$ws.Cells.Item(5, 4).Value = 82
$ws.Cells.Item(5, 6).Value = 20
$ws.Cells.Item(5, 7).Value = 0.756
$ws.Cells.Item(5, 10).Value = 80.3
$ws.Cells.Item(5, 12).Value = 6.2
$ws.Cells.Item(5, 14).Value = 0.361
$ws.Cells.Item(5, 15).Value = 18.2
$ws.Cells.Item(5, 17).Value = 0.743
$ws.Cells.Item(5, 19).Value = 32.4
$ws.Cells.Item(5, 20).Value = 44.2
$ws.Cells.Item(5, 22).Value = 14.2
$ws.Cells.Item(5, 26).Value = 20
$ws.Cells.Item(5, 28).Value = 98.59999999999999
$ws.Cells.Item(5, 29).Value = 7.3
$ws.Cells.Item(5, 35).Value = 18
$ws.Cells.Item(5, 37).Value = 13
$ws.Cells.Item(5, 38).Value = 15
$ws.Cells.Item(5, 43).Value = 26
$ws.Cells.Item(5, 49).Value = 17
$ws.Cells.Item(5, 52).Value = 9
$ws.Cells.Item(5, 58).Value = "2011-05-05"

# Row 6
$ws.Cells.Item(6, 4).Value = 82
$ws.Cells.Item(6, 5).Value = 19
$ws.Cells.Item(6, 7).Value = 0.232
$ws.Cells.Item(6, 13).Value = 18.2
$ws.Cells.Item(6, 16).Value = 25.3
$ws.Cells.Item(6, 17).Value = 0.745
$ws.Cells.Item(6, 18).Value = 10.4
$ws.Cells.Item(6, 21).Value = 21
$ws.Cells.Item(6, 23).Value = 6.6
$ws.Cells.Item(6, 29).Value = -9
$ws.Cells.Item(6, 40).Value = 23
$ws.Cells.Item(6, 42).Value = 8
$ws.Cells.Item(6, 43).Value = 24
$ws.Cells.Item(6, 46).Value = 22
$ws.Cells.Item(6, 48).Value = 18
$ws.Cells.Item(6, 49).Value = 26
$ws.Cells.Item(6, 58).Value = "2011-05-05"

# Row 7
$ws.Cells.Item(7, 4).Value = 82
$ws.Cells.Item(7, 5).Value = 57
$ws.Cells.Item(7, 7).Value = 0.695
$ws.Cells.Item(7, 10).Value = 78.8
$ws.Cells.Item(7, 11).Value = 0.475
$ws.Cells.Item(7, 14).Value = 0.365
$ws.Cells.Item(7, 15).Value = 17.5
$ws.Cells.Item(7, 17).Value = 0.777
$ws.Cells.Item(7, 22).Value = 14
$ws.Cells.Item(7, 23).Value = 6.8
$ws.Cells.Item(7, 27).Value = 20.1
$ws.Cells.Item(7, 29).Value = 4.2
$ws.Cells.Item(7, 35).Value = 13
$ws.Cells.Item(7, 41).Value = 24
$ws.Cells.Item(7, 42).Value = 27
$ws.Cells.Item(7, 45).Value = 6
$ws.Cells.Item(7, 46).Value = 14
$ws.Cells.Item(7, 47).Value = 2
$ws.Cells.Item(7, 50).Value = 24
$ws.Cells.Item(7, 58).Value = "2011-05-05"

# Row 8
$ws.Cells.Item(8, 4).Value = 82
$ws.Cells.Item(8, 5).Value = 50
$ws.Cells.Item(8, 7).Value = 0.61
$ws.Cells.Item(8, 10).Value = 80.59999999999999
$ws.Cells.Item(8, 14).Value = 0.388
$ws.Cells.Item(8, 15).Value = 22.7
$ws.Cells.Item(8, 16).Value = 29.6
$ws.Cells.Item(8, 17).Value = 0.765
$ws.Cells.Item(8, 20).Value = 42
$ws.Cells.Item(8, 26).Value = 21
$ws.Cells.Item(8, 27).Value = 22.7
$ws.Cells.Item(8, 28).Value = 107.5
$ws.Cells.Item(8, 29).Value = 4.8
$ws.Cells.Item(8, 35).Value = 6
$ws.Cells.Item(8, 44).Value = 26
$ws.Cells.Item(8, 45).Value = 5
$ws.Cells.Item(8, 47).Value = 10
$ws.Cells.Item(8, 48).Value = 16
$ws.Cells.Item(8, 58).Value = "2011-05-05"

# Row 9
$ws.Cells.Item(9, 4).Value = 82
$ws.Cells.Item(9, 6).Value = 52
$ws.Cells.Item(9, 7).Value = 0.366
$ws.Cells.Item(9, 10).Value = 81.09999999999999
$ws.Cells.Item(9, 12).Value = 5.8
$ws.Cells.Item(9, 14).Value = 0.376
$ws.Cells.Item(9, 16).Value = 22.6
$ws.Cells.Item(9, 17).Value = 0.737
$ws.Cells.Item(9, 22).Value = 13
$ws.Cells.Item(9, 25).Value = 4.6
$ws.Cells.Item(9, 34).Value = 6
$ws.Cells.Item(9, 36).Value = 13
$ws.Cells.Item(9, 37).Value = 17
$ws.Cells.Item(9, 39).Value = 22
$ws.Cells.Item(9, 42).Value = 25
$ws.Cells.Item(9, 48).Value = 2
$ws.Cells.Item(9, 49).Value = 16
$ws.Cells.Item(9, 53).Value = 27
$ws.Cells.Item(9, 58).Value = "2011-05-05"

# Row 10
$ws.Cells.Item(10, 4).Value = 82
$ws.Cells.Item(10, 5).Value = 36
$ws.Cells.Item(10, 7).Value = 0.439
$ws.Cells.Item(10, 9).Value = 39.6
$ws.Cells.Item(10, 11).Value = 0.461
$ws.Cells.Item(10, 16).Value = 20.7
$ws.Cells.Item(10, 17).Value = 0.761
$ws.Cells.Item(10, 24).Value = 5
$ws.Cells.Item(10, 26).Value = 22
$ws.Cells.Item(10, 28).Value = 103.4
$ws.Cells.Item(10, 29).Value = -2.3
$ws.Cells.Item(10, 37).Value = 14
$ws.Cells.Item(10, 52).Value = 24
$ws.Cells.Item(10, 53).Value = 30
$ws.Cells.Item(10, 58).Value = "2011-05-05"

# Row 11
$ws.Cells.Item(11, 4).Value = 82
$ws.Cells.Item(11, 6).Value = 39
$ws.Cells.Item(11, 7).Value = 0.524
$ws.Cells.Item(11, 9).Value = 38.7
$ws.Cells.Item(11, 10).Value = 85.09999999999999
$ws.Cells.Item(11, 11).Value = 0.454
$ws.Cells.Item(11, 13).Value = 22.5
$ws.Cells.Item(11, 14).Value = 0.367
$ws.Cells.Item(11, 15).Value = 20.3
$ws.Cells.Item(11, 16).Value = 25.4
$ws.Cells.Item(11, 17).Value = 0.801
$ws.Cells.Item(11, 18).Value = 11.7
$ws.Cells.Item(11, 20).Value = 42.8
$ws.Cells.Item(11, 27).Value = 21.4
$ws.Cells.Item(11, 28).Value = 105.9
$ws.Cells.Item(11, 29).Value = 2.2
$ws.Cells.Item(11, 30).Value = 1
$ws.Cells.Item(11, 31).Value = 14
$ws.Cells.Item(11, 32).Value = 14
$ws.Cells.Item(11, 33).Value = 14
$ws.Cells.Item(11, 40).Value = 9
$ws.Cells.Item(11, 42).Value = 7
$ws.Cells.Item(11, 46).Value = 7
$ws.Cells.Item(11, 47).Value = 1
$ws.Cells.Item(11, 49).Value = 21
$ws.Cells.Item(11, 52).Value = 11
$ws.Cells.Item(11, 53).Value = 8
$ws.Cells.Item(11, 58).Value = "2011-05-05"

# Row 12
$ws.Cells.Item(12, 4).Value = 82
$ws.Cells.Item(12, 6).Value = 45
$ws.Cells.Item(12, 7).Value = 0.451
$ws.Cells.Item(12, 11).Value = 0.442
$ws.Cells.Item(12, 13).Value = 20.2
$ws.Cells.Item(12, 17).Value = 0.782
$ws.Cells.Item(12, 18).Value = 11.1
$ws.Cells.Item(12, 19).Value = 32.4
$ws.Cells.Item(12, 20).Value = 43.5
$ws.Cells.Item(12, 22).Value = 15.4
$ws.Cells.Item(12, 28).Value = 99.8
$ws.Cells.Item(12, 29).Value = -1.1
$ws.Cells.Item(12, 34).Value = 20
$ws.Cells.Item(12, 48).Value = 27
$ws.Cells.Item(12, 55).Value = 19
$ws.Cells.Item(12, 58).Value = "2011-05-05"

# Row 13
$ws.Cells.Item(13, 4).Value = 82
$ws.Cells.Item(13, 6).Value = 50
$ws.Cells.Item(13, 7).Value = 0.39
$ws.Cells.Item(13, 11).Value = 0.457
$ws.Cells.Item(13, 16).Value = 26.7
$ws.Cells.Item(13, 17).Value = 0.707
$ws.Cells.Item(13, 20).Value = 42.1
$ws.Cells.Item(13, 28).Value = 98.59999999999999
$ws.Cells.Item(13, 29).Value = -3.1
$ws.Cells.Item(13, 37).Value = 20
$ws.Cells.Item(13, 41).Value = 8
$ws.Cells.Item(13, 46).Value = 9
$ws.Cells.Item(13, 50).Value = 13
$ws.Cells.Item(13, 58).Value = "2011-05-05"

# Row 14
$ws.Cells.Item(14, 4).Value = 82
$ws.Cells.Item(14, 5).Value = 57
$ws.Cells.Item(14, 7).Value = 0.695
$ws.Cells.Item(14, 10).Value = 82.40000000000001
$ws.Cells.Item(14, 11).Value = 0.463
$ws.Cells.Item(14, 14).Value = 0.352
$ws.Cells.Item(14, 15).Value = 18.8
$ws.Cells.Item(14, 17).Value = 0.779
$ws.Cells.Item(14, 18).Value = 12.1
$ws.Cells.Item(14, 24).Value = 5.1
$ws.Cells.Item(14, 26).Value = 19
$ws.Cells.Item(14, 28).Value = 101.5
$ws.Cells.Item(14, 29).Value = 6.1
$ws.Cells.Item(14, 35).Value = 10
$ws.Cells.Item(14, 42).Value = 17
$ws.Cells.Item(14, 43).Value = 7
$ws.Cells.Item(14, 44).Value = 5
$ws.Cells.Item(14, 45).Value = 7
$ws.Cells.Item(14, 48).Value = 5
$ws.Cells.Item(14, 52).Value = 2
$ws.Cells.Item(14, 53).Value = 18
$ws.Cells.Item(14, 58).Value = "2011-05-05"

# Row 15
$ws.Cells.Item(15, 4).Value = 82
$ws.Cells.Item(15, 6).Value = 36
$ws.Cells.Item(15, 7).Value = 0.5610000000000001
$ws.Cells.Item(15, 14).Value = 0.334
$ws.Cells.Item(15, 16).Value = 24.2
$ws.Cells.Item(15, 17).Value = 0.75
$ws.Cells.Item(15, 18).Value = 11.8
$ws.Cells.Item(15, 19).Value = 29.2
$ws.Cells.Item(15, 21).Value = 20.6
$ws.Cells.Item(15, 25).Value = 6.2
$ws.Cells.Item(15, 26).Value = 20.8
$ws.Cells.Item(15, 29).Value = 2.3
$ws.Cells.Item(15, 31).Value = 11
$ws.Cells.Item(15, 34).Value = 6
$ws.Cells.Item(15, 40).Value = 27
$ws.Cells.Item(15, 41).Value = 17
$ws.Cells.Item(15, 42).Value = 15
$ws.Cells.Item(15, 48).Value = 12
$ws.Cells.Item(15, 58).Value = "2011-05-05"

# Row 16
$ws.Cells.Item(16, 4).Value = 82
$ws.Cells.Item(16, 5).Value = 58
$ws.Cells.Item(16, 6).Value = 24
$ws.Cells.Item(16, 7).Value = 0.707
$ws.Cells.Item(16, 8).Value = 48.2
$ws.Cells.Item(16, 9).Value = 37
$ws.Cells.Item(16, 10).Value = 76.8
$ws.Cells.Item(16, 11).Value = 0.481
$ws.Cells.Item(16, 12).Value = 6.7
$ws.Cells.Item(16, 14).Value = 0.37
$ws.Cells.Item(16, 17).Value = 0.769
$ws.Cells.Item(16, 18).Value = 9.6
$ws.Cells.Item(16, 19).Value = 32.5
$ws.Cells.Item(16, 20).Value = 42.1
$ws.Cells.Item(16, 21).Value = 20
$ws.Cells.Item(16, 28).Value = 102.1
$ws.Cells.Item(16, 29).Value = 7.5
$ws.Cells.Item(16, 30).Value = 1
$ws.Cells.Item(16, 35).Value = 19
$ws.Cells.Item(16, 40).Value = 7
$ws.Cells.Item(16, 44).Value = 27
$ws.Cells.Item(16, 45).Value = 2
$ws.Cells.Item(16, 47).Value = 26
$ws.Cells.Item(16, 49).Value = 26
$ws.Cells.Item(16, 52).Value = 14
$ws.Cells.Item(16, 53).Value = 5
$ws.Cells.Item(16, 58).Value = "2011-05-05"

# Row 17
$ws.Cells.Item(17, 4).Value = 82
$ws.Cells.Item(17, 6).Value = 47
$ws.Cells.Item(17, 7).Value = 0.427
$ws.Cells.Item(17, 14).Value = 0.342
$ws.Cells.Item(17, 17).Value = 0.757
$ws.Cells.Item(17, 20).Value = 40.8
$ws.Cells.Item(17, 21).Value = 18.8
$ws.Cells.Item(17, 23).Value = 7.5
$ws.Cells.Item(17, 26).Value = 20.5
$ws.Cells.Item(17, 27).Value = 20.7
$ws.Cells.Item(17, 31).Value = 21
$ws.Cells.Item(17, 32).Value = 21
$ws.Cells.Item(17, 33).Value = 21
$ws.Cells.Item(17, 42).Value = 23
$ws.Cells.Item(17, 44).Value = 18
$ws.Cells.Item(17, 45).Value = 16
$ws.Cells.Item(17, 49).Value = 11
$ws.Cells.Item(17, 50).Value = 14
$ws.Cells.Item(17, 52).Value = 16
$ws.Cells.Item(17, 55).Value = 17
$ws.Cells.Item(17, 58).Value = "2011-05-05"

# Row 18
$ws.Cells.Item(18, 4).Value = 82
$ws.Cells.Item(18, 6).Value = 65
$ws.Cells.Item(18, 7).Value = 0.207
$ws.Cells.Item(18, 9).Value = 37.7
$ws.Cells.Item(18, 11).Value = 0.441
$ws.Cells.Item(18, 18).Value = 13.2
$ws.Cells.Item(18, 19).Value = 31.2
$ws.Cells.Item(18, 22).Value = 17
$ws.Cells.Item(18, 26).Value = 22.3
$ws.Cells.Item(18, 28).Value = 101.1
$ws.Cells.Item(18, 29).Value = -6.6
$ws.Cells.Item(18, 34).Value = 20
$ws.Cells.Item(18, 41).Value = 12
$ws.Cells.Item(18, 42).Value = 18
$ws.Cells.Item(18, 47).Value = 25
$ws.Cells.Item(18, 49).Value = 17
$ws.Cells.Item(18, 50).Value = 10
$ws.Cells.Item(18, 53).Value = 15
$ws.Cells.Item(18, 58).Value = "2011-05-05"

# Row 19
$ws.Cells.Item(19, 4).Value = 82
$ws.Cells.Item(19, 5).Value = 24
$ws.Cells.Item(19, 7).Value = 0.293
$ws.Cells.Item(19, 9).Value = 35.6
$ws.Cells.Item(19, 11).Value = 0.44
$ws.Cells.Item(19, 14).Value = 0.343
$ws.Cells.Item(19, 19).Value = 29.8
$ws.Cells.Item(19, 21).Value = 21
$ws.Cells.Item(19, 22).Value = 14
$ws.Cells.Item(19, 28).Value = 94.2
$ws.Cells.Item(19, 29).Value = -6.2
$ws.Cells.Item(19, 34).Value = 1
$ws.Cells.Item(19, 41).Value = 25
$ws.Cells.Item(19, 47).Value = 19
$ws.Cells.Item(19, 50).Value = 16
$ws.Cells.Item(19, 51).Value = 13
$ws.Cells.Item(19, 55).Value = 26
$ws.Cells.Item(19, 58).Value = "2011-05-05"

# Row 20
$ws.Cells.Item(20, 4).Value = 82
$ws.Cells.Item(20, 5).Value = 46
$ws.Cells.Item(20, 7).Value = 0.5610000000000001
$ws.Cells.Item(20, 14).Value = 0.36
$ws.Cells.Item(20, 15).Value = 17.7
$ws.Cells.Item(20, 18).Value = 10
$ws.Cells.Item(20, 19).Value = 30.1
$ws.Cells.Item(20, 22).Value = 13
$ws.Cells.Item(20, 23).Value = 7.6
$ws.Cells.Item(20, 25).Value = 4.8
$ws.Cells.Item(20, 26).Value = 21
$ws.Cells.Item(20, 27).Value = 20.4
$ws.Cells.Item(20, 31).Value = 11
$ws.Cells.Item(20, 32).Value = 11
$ws.Cells.Item(20, 33).Value = 11
$ws.Cells.Item(20, 35).Value = 26
$ws.Cells.Item(20, 41).Value = 23
$ws.Cells.Item(20, 42).Value = 21
$ws.Cells.Item(20, 46).Value = 24
$ws.Cells.Item(20, 47).Value = 21
$ws.Cells.Item(20, 52).Value = 18
$ws.Cells.Item(20, 54).Value = 27
$ws.Cells.Item(20, 58).Value = "2011-05-05"

# Row 21
$ws.Cells.Item(21, 4).Value = 82
$ws.Cells.Item(21, 5).Value = 42
$ws.Cells.Item(21, 7).Value = 0.512
$ws.Cells.Item(21, 11).Value = 0.457
$ws.Cells.Item(21, 12).Value = 9.300000000000001
$ws.Cells.Item(21, 15).Value = 20.6
$ws.Cells.Item(21, 16).Value = 25.5
$ws.Cells.Item(21, 19).Value = 30.1
$ws.Cells.Item(21, 20).Value = 40.5
$ws.Cells.Item(21, 21).Value = 21.4
$ws.Cells.Item(21, 23).Value = 7.6
$ws.Cells.Item(21, 24).Value = 5.8
$ws.Cells.Item(21, 26).Value = 21.3
$ws.Cells.Item(21, 28).Value = 106.5
$ws.Cells.Item(21, 29).Value = 0.8
$ws.Cells.Item(21, 34).Value = 20
$ws.Cells.Item(21, 37).Value = 19
$ws.Cells.Item(21, 40).Value = 8
$ws.Cells.Item(21, 42).Value = 6
$ws.Cells.Item(21, 58).Value = "2011-05-05"

# Row 22
$ws.Cells.Item(22, 4).Value = 82
$ws.Cells.Item(22, 5).Value = 55
$ws.Cells.Item(22, 7).Value = 0.671
$ws.Cells.Item(22, 9).Value = 37.4
$ws.Cells.Item(22, 11).Value = 0.464
$ws.Cells.Item(22, 12).Value = 5.9
$ws.Cells.Item(22, 14).Value = 0.347
$ws.Cells.Item(22, 15).Value = 24.1
$ws.Cells.Item(22, 16).Value = 29.3
$ws.Cells.Item(22, 17).Value = 0.823
$ws.Cells.Item(22, 19).Value = 31.8
$ws.Cells.Item(22, 23).Value = 8
$ws.Cells.Item(22, 26).Value = 22.4
$ws.Cells.Item(22, 29).Value = 3.8
$ws.Cells.Item(22, 34).Value = 1
$ws.Cells.Item(22, 40).Value = 19
$ws.Cells.Item(22, 45).Value = 8
$ws.Cells.Item(22, 46).Value = 8
$ws.Cells.Item(22, 47).Value = 24
$ws.Cells.Item(22, 48).Value = 15
$ws.Cells.Item(22, 50).Value = 2
$ws.Cells.Item(22, 51).Value = 7
$ws.Cells.Item(22, 58).Value = "2011-05-05"

# Row 23
$ws.Cells.Item(23, 4).Value = 82
$ws.Cells.Item(23, 5).Value = 52
$ws.Cells.Item(23, 7).Value = 0.634
$ws.Cells.Item(23, 11).Value = 0.461
$ws.Cells.Item(23, 13).Value = 25.6
$ws.Cells.Item(23, 14).Value = 0.366
$ws.Cells.Item(23, 16).Value = 25.6
$ws.Cells.Item(23, 17).Value = 0.6919999999999999
$ws.Cells.Item(23, 19).Value = 32.7
$ws.Cells.Item(23, 20).Value = 43.2
$ws.Cells.Item(23, 28).Value = 99.2
$ws.Cells.Item(23, 29).Value = 5.5
$ws.Cells.Item(23, 37).Value = 16
$ws.Cells.Item(23, 41).Value = 22
$ws.Cells.Item(23, 42).Value = 5
$ws.Cells.Item(23, 44).Value = 17
$ws.Cells.Item(23, 45).Value = 1
$ws.Cells.Item(23, 53).Value = 4
$ws.Cells.Item(23, 54).Value = 16
$ws.Cells.Item(23, 55).Value = 5
$ws.Cells.Item(23, 58).Value = "2011-05-05"

# Row 24
$ws.Cells.Item(24, 4).Value = 82
$ws.Cells.Item(24, 6).Value = 41
$ws.Cells.Item(24, 7).Value = 0.5
$ws.Cells.Item(24, 11).Value = 0.461
$ws.Cells.Item(24, 14).Value = 0.355
$ws.Cells.Item(24, 15).Value = 17.4
$ws.Cells.Item(24, 16).Value = 22.6
$ws.Cells.Item(24, 17).Value = 0.77
$ws.Cells.Item(24, 19).Value = 31.4
$ws.Cells.Item(24, 20).Value = 41.8
$ws.Cells.Item(24, 22).Value = 13
$ws.Cells.Item(24, 23).Value = 7.6
$ws.Cells.Item(24, 28).Value = 99
$ws.Cells.Item(24, 29).Value = 1.5
$ws.Cells.Item(24, 31).Value = 16
$ws.Cells.Item(24, 32).Value = 16
$ws.Cells.Item(24, 33).Value = 16
$ws.Cells.Item(24, 35).Value = 11
$ws.Cells.Item(24, 37).Value = 15
$ws.Cells.Item(24, 41).Value = 26
$ws.Cells.Item(24, 42).Value = 26
$ws.Cells.Item(24, 43).Value = 10
$ws.Cells.Item(24, 49).Value = 10
$ws.Cells.Item(24, 51).Value = 14
$ws.Cells.Item(24, 53).Value = 28
$ws.Cells.Item(24, 54).Value = 18
$ws.Cells.Item(24, 55).Value = 13
$ws.Cells.Item(24, 58).Value = "2011-05-05"

# Row 25
$ws.Cells.Item(25, 4).Value = 82
$ws.Cells.Item(25, 5).Value = 40
$ws.Cells.Item(25, 6).Value = 42
$ws.Cells.Item(25, 10).Value = 83.5
$ws.Cells.Item(25, 13).Value = 22.6
$ws.Cells.Item(25, 14).Value = 0.377
$ws.Cells.Item(25, 15).Value = 18
$ws.Cells.Item(25, 17).Value = 0.759
$ws.Cells.Item(25, 19).Value = 30.2
$ws.Cells.Item(25, 20).Value = 40.2
$ws.Cells.Item(25, 21).Value = 23.7
$ws.Cells.Item(25, 22).Value = 14.3
$ws.Cells.Item(25, 25).Value = 4.3
$ws.Cells.Item(25, 27).Value = 21.2
$ws.Cells.Item(25, 28).Value = 105
$ws.Cells.Item(25, 30).Value = 1
$ws.Cells.Item(25, 45).Value = 17
$ws.Cells.Item(25, 46).Value = 23
$ws.Cells.Item(25, 47).Value = 3
$ws.Cells.Item(25, 48).Value = 19
$ws.Cells.Item(25, 49).Value = 25
$ws.Cells.Item(25, 51).Value = 7
$ws.Cells.Item(25, 53).Value = 12
$ws.Cells.Item(25, 58).Value = "2011-05-05"

# Row 26
$ws.Cells.Item(26, 4).Value = 82
$ws.Cells.Item(26, 5).Value = 48
$ws.Cells.Item(26, 7).Value = 0.585
$ws.Cells.Item(26, 9).Value = 36
$ws.Cells.Item(26, 10).Value = 80.5
$ws.Cells.Item(26, 14).Value = 0.345
$ws.Cells.Item(26, 15).Value = 18
$ws.Cells.Item(26, 16).Value = 22.4
$ws.Cells.Item(26, 17).Value = 0.804
$ws.Cells.Item(26, 18).Value = 12.1
$ws.Cells.Item(26, 19).Value = 27.2
$ws.Cells.Item(26, 20).Value = 39.3
$ws.Cells.Item(26, 21).Value = 21.2
$ws.Cells.Item(26, 28).Value = 96.3
$ws.Cells.Item(26, 29).Value = 1.5
$ws.Cells.Item(26, 30).Value = 1
$ws.Cells.Item(26, 35).Value = 25
$ws.Cells.Item(26, 36).Value = 19
$ws.Cells.Item(26, 38).Value = 13
$ws.Cells.Item(26, 40).Value = 21
$ws.Cells.Item(26, 41).Value = 18
$ws.Cells.Item(26, 44).Value = 4
$ws.Cells.Item(26, 46).Value = 27
$ws.Cells.Item(26, 48).Value = 4
$ws.Cells.Item(26, 53).Value = 11
$ws.Cells.Item(26, 55).Value = 12
$ws.Cells.Item(26, 58).Value = "2011-05-05"

# Row 27
$ws.Cells.Item(27, 4).Value = 82
$ws.Cells.Item(27, 5).Value = 24
$ws.Cells.Item(27, 7).Value = 0.293
$ws.Cells.Item(27, 15).Value = 17.7
$ws.Cells.Item(27, 16).Value = 24.2
$ws.Cells.Item(27, 18).Value = 13.1
$ws.Cells.Item(27, 22).Value = 16.1
$ws.Cells.Item(27, 25).Value = 5.7
$ws.Cells.Item(27, 27).Value = 21
$ws.Cells.Item(27, 28).Value = 99.40000000000001
$ws.Cells.Item(27, 29).Value = -5.3
$ws.Cells.Item(27, 41).Value = 21
$ws.Cells.Item(27, 42).Value = 15
$ws.Cells.Item(27, 47).Value = 23
$ws.Cells.Item(27, 52).Value = 23
$ws.Cells.Item(27, 53).Value = 14
$ws.Cells.Item(27, 58).Value = "2011-05-05"

# Row 28
$ws.Cells.Item(28, 4).Value = 82
$ws.Cells.Item(28, 5).Value = 61
$ws.Cells.Item(28, 7).Value = 0.744
$ws.Cells.Item(28, 9).Value = 38.4
$ws.Cells.Item(28, 12).Value = 8.4
$ws.Cells.Item(28, 14).Value = 0.397
$ws.Cells.Item(28, 15).Value = 18.5
$ws.Cells.Item(28, 16).Value = 24.2
$ws.Cells.Item(28, 17).Value = 0.767
$ws.Cells.Item(28, 19).Value = 31.7
$ws.Cells.Item(28, 21).Value = 22.4
$ws.Cells.Item(28, 23).Value = 7.3
$ws.Cells.Item(28, 25).Value = 4.6
$ws.Cells.Item(28, 28).Value = 103.7
$ws.Cells.Item(28, 29).Value = 5.7
$ws.Cells.Item(28, 34).Value = 20
$ws.Cells.Item(28, 35).Value = 5
$ws.Cells.Item(28, 38).Value = 4
$ws.Cells.Item(28, 41).Value = 11
$ws.Cells.Item(28, 42).Value = 14
$ws.Cells.Item(28, 45).Value = 9
$ws.Cells.Item(28, 46).Value = 12
$ws.Cells.Item(28, 49).Value = 14
$ws.Cells.Item(28, 58).Value = "2011-05-05"

# Row 29
$ws.Cells.Item(29, 4).Value = 82
$ws.Cells.Item(29, 6).Value = 60
$ws.Cells.Item(29, 7).Value = 0.268
$ws.Cells.Item(29, 9).Value = 38.3
$ws.Cells.Item(29, 10).Value = 82.40000000000001
$ws.Cells.Item(29, 11).Value = 0.465
$ws.Cells.Item(29, 19).Value = 28.6
$ws.Cells.Item(29, 25).Value = 5.6
$ws.Cells.Item(29, 26).Value = 22
$ws.Cells.Item(29, 27).Value = 19.8
$ws.Cells.Item(29, 34).Value = 20
$ws.Cells.Item(29, 35).Value = 7
$ws.Cells.Item(29, 36).Value = 12
$ws.Cells.Item(29, 42).Value = 19
$ws.Cells.Item(29, 43).Value = 22
$ws.Cells.Item(29, 46).Value = 21
$ws.Cells.Item(29, 49).Value = 21
$ws.Cells.Item(29, 52).Value = 25
$ws.Cells.Item(29, 55).Value = 27
$ws.Cells.Item(29, 58).Value = "2011-05-05"

# Row 30
$ws.Cells.Item(30, 4).Value = 82
$ws.Cells.Item(30, 6).Value = 43
$ws.Cells.Item(30, 7).Value = 0.476
$ws.Cells.Item(30, 9).Value = 37.4
$ws.Cells.Item(30, 10).Value = 80.40000000000001
$ws.Cells.Item(30, 11).Value = 0.465
$ws.Cells.Item(30, 12).Value = 5.3
$ws.Cells.Item(30, 13).Value = 15.3
$ws.Cells.Item(30, 14).Value = 0.346
$ws.Cells.Item(30, 15).Value = 19.4
$ws.Cells.Item(30, 16).Value = 25.1
$ws.Cells.Item(30, 17).Value = 0.771
$ws.Cells.Item(30, 18).Value = 11
$ws.Cells.Item(30, 19).Value = 28.5
$ws.Cells.Item(30, 20).Value = 39.5
$ws.Cells.Item(30, 21).Value = 23.4
$ws.Cells.Item(30, 22).Value = 14.3
$ws.Cells.Item(30, 25).Value = 5
$ws.Cells.Item(30, 26).Value = 22.7
$ws.Cells.Item(30, 27).Value = 22
$ws.Cells.Item(30, 28).Value = 99.40000000000001
$ws.Cells.Item(30, 29).Value = -1.8
$ws.Cells.Item(30, 30).Value = 1
$ws.Cells.Item(30, 31).Value = 18
$ws.Cells.Item(30, 32).Value = 18
$ws.Cells.Item(30, 33).Value = 18
$ws.Cells.Item(30, 35).Value = 15
$ws.Cells.Item(30, 36).Value = 21
$ws.Cells.Item(30, 40).Value = 20
$ws.Cells.Item(30, 50).Value = 3
$ws.Cells.Item(30, 51).Value = 19
$ws.Cells.Item(30, 53).Value = 3
$ws.Cells.Item(30, 58).Value = "2011-05-05"

# Row 31
$ws.Cells.Item(31, 4).Value = 82
$ws.Cells.Item(31, 6).Value = 59
$ws.Cells.Item(31, 7).Value = 0.28
$ws.Cells.Item(31, 9).Value = 37.2
$ws.Cells.Item(31, 10).Value = 84
$ws.Cells.Item(31, 14).Value = 0.332
$ws.Cells.Item(31, 15).Value = 18.2
$ws.Cells.Item(31, 16).Value = 24.4
$ws.Cells.Item(31, 17).Value = 0.745
$ws.Cells.Item(31, 18).Value = 12.4
$ws.Cells.Item(31, 20).Value = 41.3
$ws.Cells.Item(31, 22).Value = 15.3
$ws.Cells.Item(31, 24).Value = 6.1
$ws.Cells.Item(31, 25).Value = 5
$ws.Cells.Item(31, 26).Value = 22.6
$ws.Cells.Item(31, 27).Value = 20.3
$ws.Cells.Item(31, 28).Value = 97.3
$ws.Cells.Item(31, 29).Value = -7.4
$ws.Cells.Item(31, 31).Value = 27
$ws.Cells.Item(31, 32).Value = 27
$ws.Cells.Item(31, 33).Value = 27
$ws.Cells.Item(31, 35).Value = 17
$ws.Cells.Item(31, 38).Value = 28
$ws.Cells.Item(31, 40).Value = 28
$ws.Cells.Item(31, 42).Value = 12
$ws.Cells.Item(31, 43).Value = 25
$ws.Cells.Item(31, 46).Value = 15
$ws.Cells.Item(31, 48).Value = 26
$ws.Cells.Item(31, 51).Value = 20
$ws.Cells.Item(31, 58).Value = "2011-05-05"

